# V. 75 "A descubierto"
# Adds a new movie entry ("A descubierto") to the ranked table on the
# "Películas" sheet. The table is sorted descending by "Puntuación total"
# (column C), and the new entry's score slots it in at row 96, pushing the
# previous rows 96-98 ("Alimañas", "Bajo el mismo techo", "Chicos buenos")
# down to rows 97-99.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")

# Insert a blank row at row 96; existing rows 96-98 shift down to 97-99
# (formulas and row-relative references are adjusted automatically).
$null = $ws.Rows("96:96").Insert()

# Populate the new entry.
$ws.Range("B96").Value = "A descubierto"
$ws.Range("C96").Formula = "=AVERAGE(D96,E96,E96,F96,G96,H96,H96,I96)"
$ws.Range("D96").Value = 3
$ws.Range("E96").Value = 4
$ws.Range("F96").Value = 4
$ws.Range("G96").Value = 4
$ws.Range("H96").Value = 5.4
$ws.Range("I96").Value = 4.7

# Grow the table ("Tabla24") so the new row becomes part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:I99"))

# The previously-newest entry ("Garra", row 40) loses its "latest addition"
# highlight now that "A descubierto" is the newest entry.
$ws.Range("B40").HorizontalAlignment = -4131

# Update the saved view/selection state.
$null = $ws.Range("C85").Select()
$excel.ActiveWindow.ScrollRow = 77
